$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$startRow = 2
$endRow = 113
$col = "C"

for ($row = $startRow; $row -le $endRow; $row++) {
    $ws.Range("$col$row").Value = 46061
}
